# ---------------------------------------------------------------------------
# Applies the "improved recommendation on homepage" edits:
#  1. De-emphasize (remove red highlight) on the word "audio" in two spots.
#  2. Replace "the web app" -> "our server" in the video-upload pipeline
#     paragraph (5 occurrences), and drop a redundant "the uploading" -> "uploading".
#  3. Shorten "for the purpose of" -> "for" (also relocates the trailing
#     "_GoBack" bookmark to that edit point, matching Word's own behaviour of
#     anchoring _GoBack at the most recent edit).
#  4. Swap "Azure Web App" for "Azure Virtual Machine" in the technologies
#     paragraph and reword the hosting sentence accordingly.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Black / theme "text1" color sentinel -- matches the existing black runs
# that already sit next to the red "audio" runs in this document.
$blackText1 = -587137025

# ---------------------------------------------------------------------------
# 1. Remove the red color on "audio" (two occurrences) by replacing the
#    surrounding phrase with explicitly-formatted (black) replacement text.
#    Using Find.Replacement.Font + Format:=$true both re-colors and merges
#    the matched runs into a single, uniformly-formatted run.
# ---------------------------------------------------------------------------

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = " in audio file"
$find.Replacement.Text = " in audio file"
$find.Replacement.Font.Color = $blackText1
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2, $true, $false, $false, $false) | Out-Null

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Text = "is audio or visual. Furthermore, the service "
$find2.Replacement.Text = "is audio or visual. Furthermore, the service "
$find2.Replacement.Font.Color = $blackText1
$find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2, $true, $false, $false, $false) | Out-Null

# ---------------------------------------------------------------------------
# 2. "the web app" -> "our server" across the upload-pipeline paragraph.
#    These runs already share uniform (non-red) formatting, so plain text
#    replacement is sufficient.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("all its information will first be passed to the web app, then stored", $true, $false, $false, $false, $false, $true, 1, $false, "all its information will first be passed to our server, then stored", 2) | Out-Null

$d.Content.Find.Execute("After each chunk is been uploaded to the web app, it will", $true, $false, $false, $false, $false, $true, 1, $false, "After each chunk is been uploaded to our server, it will", 2) | Out-Null

$d.Content.Find.Execute("a special post request will be sent to the web app signaling the end of the uploading process", $true, $false, $false, $false, $false, $true, 1, $false, "a special post request will be sent to our server signaling the end of uploading process", 2) | Out-Null

$d.Content.Find.Execute("Finally, the web app will perform", $true, $false, $false, $false, $false, $true, 1, $false, "Finally, our server will perform", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. "for the purpose of" -> "for", and relocate the "_GoBack" bookmark to
#    that exact edit point (mirrors where Word itself drops _GoBack: the
#    location of the user's most recent edit).
# ---------------------------------------------------------------------------

$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$rngShrink = $d.Content
$rngShrink.Find.Execute("cut in slices for the purpose of preventing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngShrink.Text = "cut in slices for preventing"

$rngAnchor = $d.Content
$rngAnchor.Find.Execute("cut in slices ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPoint = $d.Range($rngAnchor.End, $rngAnchor.End)
$d.Bookmarks.Add("_GoBack", $anchorPoint) | Out-Null

# ---------------------------------------------------------------------------
# 4. Azure Web App -> Azure Virtual Machine (technologies paragraph).
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("We are using three key services from Azure: Azure Web App, Azure SQL", $true, $false, $false, $false, $false, $true, 1, $false, "We are using three key services from Azure: Azure Virtual Machine, Azure SQL", 2) | Out-Null

$d.Content.Find.Execute("Our hosting service is the Azure Web App running on Python version 2.7.13, where Azure SQL", $true, $false, $false, $false, $false, $true, 1, $false, "We are hosting our service on Azure Virtual Machine, where Azure SQL", 2) | Out-Null

Write-Host "done"
